$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.202.44'
$ws.Range("E2").Value = '  +0.94%  '

$ws.Range("D3").Value = '1.798.45'
$ws.Range("E3").Value = '  +2.35%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.84'
$ws.Range("E5").Value = '  +0.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4749'
$ws.Range("E7").Value = '  +25.86%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3721'
$ws.Range("E8").Value = '  +11.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.55'
$ws.Range("E9").Value = '  -0.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07714'
$ws.Range("E10").Value = '  +7.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.145'
$ws.Range("E11").Value = '  +2.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.69'
$ws.Range("E12").Value = '  +2.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.002'
$ws.Range("E13").Value = '  -0.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.344'
$ws.Range("E14").Value = '  +2.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.342'
$ws.Range("E15").Value = '  +2.45%  '

$ws.Range("D16").Value = '1.796.15'
$ws.Range("E16").Value = '  +2.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001097'
$ws.Range("E17").Value = '  +4.34%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06741'
$ws.Range("E18").Value = '  +2.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.02'
$ws.Range("E19").Value = '  +2.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.37'
$ws.Range("E21").Value = '  +2.79%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.423'

$ws.Range("D23").Value = '28.199.60'
$ws.Range("E23").Value = '  +0.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.00'
$ws.Range("E24").Value = '  +2.69%  '

$ws.Range("E25").Value = '  +1.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.68'
$ws.Range("E26").Value = '  +4.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.409'
$ws.Range("E27").Value = '  +3.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.71'
$ws.Range("E28").Value = '  -1.00%  '

$ws.Range("D29").Value = '2.002.63'
$ws.Range("E29").Value = '  +2.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.40'
$ws.Range("E30").Value = '  +1.93%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.277'
$ws.Range("E31").Value = '  +0.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.042'
$ws.Range("E32").Value = '  +0.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09662'
$ws.Range("E33").Value = '  +10.56%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.949'
$ws.Range("E34").Value = '  +3.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02379'
$ws.Range("E35").Value = '  +1.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.21'
$ws.Range("E36").Value = '  +0.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6706'
$ws.Range("E37").Value = '  +2.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06319'
$ws.Range("E38").Value = '  +2.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.249'
$ws.Range("E39").Value = '  +2.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2191'
$ws.Range("E40").Value = '  +4.30%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.222'
$ws.Range("E41").Value = '  +1.06%  '

$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.483'
$ws.Range("E42").Value = '  +2.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.129'
$ws.Range("E43").Value = '  +1.64%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.12'
$ws.Range("E44").Value = '  +2.64%  '

$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6182'
$ws.Range("E46").Value = '  +2.64%  '

$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.875'
$ws.Range("E47").Value = '  +1.35%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.02'
$ws.Range("E48").Value = '  -0.89%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.053'
$ws.Range("E49").Value = '  +2.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.176'
$ws.Range("E50").Value = '  -0.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07098'
$ws.Range("E51").Value = '  -0.75%  '
